# Auto-update epexspot_prices.xlsx with the latest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Prix Spot": append a new date column (CG) with hourly prices.
# ---------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("CG1").Value = "06-sep"
$wsPrix.Range("CF1").Copy()
$wsPrix.Range("CG1").PasteSpecial(-4122)  # xlPasteFormats, reuse the header style

$prixValues = @(74.64, 57.8, 41.35, 20.66, 20, 25.38, 25.49, 25.86, 29.28, 15.23, 1.72, 0, -0.02, -0.86, -0.99, -0.01, -0.01, 0, 17.42, 24.37, 42, 19.01, 36.26, 27.99)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 85).Value = $prixValues[$i]
}

# ---------------------------------------------------------------
# Sheet "Gaz": append a new daily row.
# ---------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date to be stored as plain text (matching the other A-column
# cells) instead of letting Excel auto-convert it to a date serial number.
$gazDate = "2025-09-04"
$wsGaz.Range("A82").NumberFormat = "@"
$wsGaz.Range("A82").Value = $gazDate

$wsGaz.Range("B82").Value = 31.5

# ---------------------------------------------------------------
# Sheet "CO2": append a new daily row.
# ---------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$co2Date = "2025-09-04"
$wsCo2.Range("A82").NumberFormat = "@"
$wsCo2.Range("A82").Value = $co2Date

$wsCo2.Range("B82").Value = 74.9
